$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B79").Value = 45300
$ws.Range("C79").Value = "4 botellones"
$ws.Range("D79").Value = -212

[void]$ws.Range("D80").Select()
